$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6's date of birth was mistakenly entered with the wrong year (2024 instead of 2004),
# which made the "can vote" age check throw the wrong error popup. Fix the value.
$ws.Range("I6").Value = "6/20/2004"

# Row 7's date of birth was stored as plain text ("7/30/2003"), which also broke the
# vote-eligibility check and used the wrong year. Give it the same date formatting as the
# other Date of Birth cells, then set it to a proper date value.
$ws.Range("I6").Copy()
$ws.Range("I7").PasteSpecial(-4122)
$ws.Range("I7").Value = "7/30/2004"

$excel.CutCopyMode = 0

# Leave the selection where it ended up after editing the two date cells.
$ws.Range("I8").Select()

$wb.Save()
